$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename customer "Mr Mustacheo" -> "MR. MUSTACHEO GROUP LLC" (row 25, column A)
$ws.Range("A25").Value = "MR. MUSTACHEO GROUP LLC"

# Update salesperson code for row 25 (column C) from "013" to "023"
$ws.Range("C25").Value = "023"

# Set "Last Invoice Date" for row 25 (column D) to 2025-09-11 (serial 45911),
# copying the date format from a neighboring already-formatted date cell (D24)
# so the style/number-format matches exactly instead of creating a new style.
$ws.Range("D25").Value = 45911
$ws.Range("D24").Copy() | Out-Null
$ws.Range("D25").PasteSpecial(-4122) | Out-Null

# Set "Last Invoice Date" for row 31 (column D) to 2025-09-11 (serial 45911) as well
$ws.Range("D31").Value = 45911
$ws.Range("D30").Copy() | Out-Null
$ws.Range("D31").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0
